# mergesort_animation.pptx — renumber a block of MERGE_SORT/MERGE labels
# (indices (4,x)/(x,7) -> (5,x)/(x,8), etc.) on slide 2, and tidy up a
# previously-split run on slide 1 back into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1: "Rectangle 29" — collapse the two runs
#   "MERGE_SORT(4, 4) " + "-> [3]"
# back into a single run "MERGE_SORT(4, 4) -> [3]" (no text change, just
# a structural re-typing of the whole box which PowerPoint coalesces
# into one run).
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$shp = $slide1.Shapes.Item("Rectangle 29")
$tr = $shp.TextFrame.TextRange
# Force a real change first so PowerPoint actually re-flows the runs,
# then set the final text, which collapses everything into one run.
$tr.Text = "~"
$tr.Text = "MERGE_SORT(4, 4) -> [3]"

# ---------------------------------------------------------------------
# Slide 2: bump each of the (start, end) pairs by +1 by retyping just
# the leading "MERGE_SORT(a, b) " / "MERGE(a, b) " portion of each
# label, leaving the trailing "-> [...]" portion (and any leading
# space) as its own, untouched run(s).
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)

function Set-LabelPrefix($ShapeName, $Start, $Length, $NewPrefix) {
    $shape = $slide2.Shapes.Item($ShapeName)
    $range = $shape.TextFrame.TextRange
    $sub = $range.Characters($Start, $Length)
    $sub.Text = $NewPrefix
}

# Rectangle 1:  "MERGE_SORT(4, 5) " -> "MERGE_SORT(5, 6) "
Set-LabelPrefix "Rectangle 1" 1 17 "MERGE_SORT(5, 6) "

# Rectangle 8:  "MERGE_SORT(4, 4) " -> "MERGE_SORT(5, 5) "
Set-LabelPrefix "Rectangle 8" 1 17 "MERGE_SORT(5, 5) "

# Rectangle 16: "MERGE_SORT(5, 5) " -> "MERGE_SORT(6, 6) "
Set-LabelPrefix "Rectangle 16" 1 17 "MERGE_SORT(6, 6) "

# Rectangle 23: " MERGE(4, 5) " -> " MERGE(5, 6) " (leading space kept as-is)
Set-LabelPrefix "Rectangle 23" 2 12 "MERGE(5, 6) "

# Rectangle 32: "MERGE_SORT(6, 7) " -> "MERGE_SORT(7, 8) "
Set-LabelPrefix "Rectangle 32" 1 17 "MERGE_SORT(7, 8) "

# Rectangle 36: " MERGE_SORT(6, 6) " -> " MERGE_SORT(7, 7) " (leading space kept as-is)
Set-LabelPrefix "Rectangle 36" 2 17 "MERGE_SORT(7, 7) "

# Rectangle 46: " MERGE_SORT(7, 7) " -> " MERGE_SORT(8, 8) " (leading space kept as-is)
Set-LabelPrefix "Rectangle 46" 2 17 "MERGE_SORT(8, 8) "

# Rectangle 51: " MERGE(6, 7) " -> " MERGE(7, 8) " (leading space kept as-is)
Set-LabelPrefix "Rectangle 51" 2 12 "MERGE(7, 8) "

# Rectangle 59: "MERGE(4, 7) " -> "MERGE(5, 8) "
Set-LabelPrefix "Rectangle 59" 1 12 "MERGE(5, 8) "

Write-Output "done"
